# Update cryptos.xlsx price/volume data (GitHub Actions scrape refresh).
# Generated from the authoritative cell-by-cell OOXML diff: for each changed
# cell we set the literal display text. Price cells whose new text parses as
# a plain decimal (single '.') are forced to stay text (NumberFormat "@")
# before assignment, matching the sheet's existing inline-string/text storage
# (otherwise COM Value assignment would silently coerce them to numbers and
# drop things like trailing zeros, e.g. "325.00" -> 325).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.733.97'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '2.434.28'
$ws.Range("E3").Value = '  +9.24%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.00'
$ws.Range("E5").Value = '  +10.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.75'
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("E7").Value = '  +2.86%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.657'
$ws.Range("E9").Value = '  +8.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.99'
$ws.Range("E10").Value = '  -2.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0948'
$ws.Range("E11").Value = '  +3.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.90'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.05'
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.35'
$ws.Range("E14").Value = '  +15.64%  '
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").Value = '2.801.31'
$ws.Range("E16").Value = '  +9.42%  '
$ws.Range("D17").Value = '2.433.95'
$ws.Range("E17").Value = '  +9.18%  '
$ws.Range("D18").Value = '43.783.33'
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("E19").Value = '  +4.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.41'
$ws.Range("E20").Value = '  +2.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.27'
$ws.Range("E21").Value = '  +3.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.52'
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '262.10'
$ws.Range("E23").Value = '  +13.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.50'
$ws.Range("E24").Value = '  +5.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.57'
$ws.Range("E25").Value = '  +4.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.18'
$ws.Range("E26").Value = '  +5.12%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.83'
$ws.Range("E28").Value = '  +3.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '23.12'
$ws.Range("E29").Value = '  +10.14%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.25'
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '177.23'
$ws.Range("E31").Value = '  +1.82%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0939'
$ws.Range("E33").Value = '  +5.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.03'
$ws.Range("E34").Value = '  +6.15%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("E35").Value = '  +5.20%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.99'
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.16'
$ws.Range("E37").Value = '  -2.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0376'
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.109'
$ws.Range("E39").Value = '  +2.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("E40").Value = '  +20.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.58'
$ws.Range("E41").Value = '  +19.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.96'
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.234'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '123.96'
$ws.Range("E44").Value = '  +21.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.03'
$ws.Range("E45").Value = '  +3.28%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.77'
$ws.Range("E47").Value = '  +6.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.47'
$ws.Range("E48").Value = '  +10.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.33'
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.587.43'
$ws.Range("E50").Value = '  +11.91%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.101'
$ws.Range("E51").Value = '  +3.83%  '
